$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: return a Range starting right after the first occurrence of
# $anchorText (to the end of the document). Used to disambiguate which
# occurrence of a repeated label (e.g. "Purpose", "Features") to operate on.
# ---------------------------------------------------------------------------
function Get-RangeAfter($doc, $anchorText) {
    $r = $doc.Content
    $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r.Collapse(0)
    $r.End = $doc.Content.End
    return $r
}

# ---------------------------------------------------------------------------
# Helper: within $searchRange, find $targetText and replace the run(s) that
# hold it with a single freshly built run described by $newRunInnerXml
# (e.g. "<w:r><w:rPr>...</w:rPr><w:t>...</w:t></w:r>"). Works by inserting
# the new run immediately before the match and then blanking out the
# (now shifted) original text run-span.
# ---------------------------------------------------------------------------
function Replace-RunXml($doc, $searchRange, $targetText, $newRunInnerXml) {
    $searchRange.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $s = $searchRange.Start
    $e = $searchRange.End
    $collapsed = $doc.Range($s, $s)
    $pkgXml = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p>" + $newRunInnerXml + "</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
    $collapsed.InsertXML($pkgXml)
    $len = $e - $s
    $oldRng = $doc.Range($s + $len, $e + $len)
    $oldRng.Text = ""
}

# ---------------------------------------------------------------------------
# 1. Insert a new "Update 6" paragraph right before
#    "Apps That Are Needed (Gap-Filling Ideas)" (after the existing
#    "Update 3" paragraph), matching that paragraph's bold/sz formatting.
# ---------------------------------------------------------------------------
$update3Rng = $d.Content
$update3Rng.Find.Execute("Update 3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$update3Para = $update3Rng.Paragraphs(1)
$update3Para.Range.InsertParagraphAfter()
$newUpdatePara = $update3Para.Next()
$newUpdatePara.Range.InsertBefore("Update 6")

# ---------------------------------------------------------------------------
# 2. Merge the split runs (removing the w:proofErr gramStart/gramEnd markers
#    along the way) back into single runs with identical visible text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(": Consolidate administrative, training, and operational data for easy access.", $true, $false, $false, $false, $false, $true, 1, $false, ": Consolidate administrative, training, and operational data for easy access.", 2) | Out-Null

$d.Content.Find.Execute("Task tracking categorized by soldier, equipment, or mission.", $true, $false, $false, $false, $false, $true, 1, $false, "Task tracking categorized by soldier, equipment, or mission.", 2) | Out-Null

$d.Content.Find.Execute(": Simplify tracking vehicle, weapon, and supply readiness at the platoon level.", $true, $false, $false, $false, $false, $true, 1, $false, ": Simplify tracking vehicle, weapon, and supply readiness at the platoon level.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Relocate the w:lastRenderedPageBreak markers:
#    a) off "Purpose" (section 3, Maintenance and Supply Tracker) and onto
#       the "Digital checklist..." bullet right below it.
#    b) off "6. Time Management and Scheduling App" and onto the "Features"
#       bullet of that same section.
#    c) off "9. Team Climate and Morale Monitor" and onto the "Features"
#       bullet of that same section.
# ---------------------------------------------------------------------------
$rngPurpose = Get-RangeAfter $d "3. Maintenance and Supply Tracker"
Replace-RunXml $d $rngPurpose "Purpose" "<w:r w:rsidRPr='008E2B9A'><w:rPr><w:b/><w:bCs/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:t>Purpose</w:t></w:r>"

Replace-RunXml $d $d.Content "Digital checklist for scheduled maintenance (linked to GCSS-Army)." "<w:r w:rsidRPr='008E2B9A'><w:rPr><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:lastRenderedPageBreak/><w:t>Digital checklist for scheduled maintenance (linked to GCSS-Army).</w:t></w:r>"

Replace-RunXml $d $d.Content "6. Time Management and Scheduling App" "<w:r w:rsidRPr='008E2B9A'><w:rPr><w:b/><w:bCs/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:t>6. Time Management and Scheduling App</w:t></w:r>"

$rngFeatures6 = Get-RangeAfter $d "6. Time Management and Scheduling App"
Replace-RunXml $d $rngFeatures6 "Features" "<w:r w:rsidRPr='008E2B9A'><w:rPr><w:b/><w:bCs/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:lastRenderedPageBreak/><w:t>Features</w:t></w:r>"

Replace-RunXml $d $d.Content "9. Team Climate and Morale Monitor" "<w:r w:rsidRPr='008E2B9A'><w:rPr><w:b/><w:bCs/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:t>9. Team Climate and Morale Monitor</w:t></w:r>"

$rngFeatures9 = Get-RangeAfter $d "9. Team Climate and Morale Monitor"
Replace-RunXml $d $rngFeatures9 "Features" "<w:r w:rsidRPr='008E2B9A'><w:rPr><w:b/><w:bCs/><w:sz w:val='22'/><w:szCs w:val='22'/></w:rPr><w:lastRenderedPageBreak/><w:t>Features</w:t></w:r>"

Write-Output "applied edits"
